# Update KRA page to enable review option based on reviewCycle
#
# This script updates the "KeyResult Data" export sheet:
#  - Column I (windowId -> here used as a "reviewable"/count flag) gets a
#    value of 1 for every existing data row that was previously blank.
#  - A few rating/count values are corrected for rows 6-11 (H column),
#    D9:D11 weight-ish values, and F7.
#  - A new trailing row (13) is appended for a window that has no
#    ratings/description yet (reviewCycle placeholder row), mirroring the
#    empty-string cells the exporter emits for "no value yet" text fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---------------------------------------------------------------
$ws.Cells.Item(4, 9).Value = 1      # I4

# --- Row 5 ---------------------------------------------------------------
$ws.Cells.Item(5, 9).Value = 1      # I5

# --- Row 6 ---------------------------------------------------------------
$ws.Cells.Item(6, 8).Value = 2      # H6
$ws.Cells.Item(6, 9).Value = 1      # I6

# --- Row 7 ---------------------------------------------------------------
$ws.Cells.Item(7, 6).Value = 20     # F7
$ws.Cells.Item(7, 8).Value = 2      # H7
$ws.Cells.Item(7, 9).Value = 1      # I7

# --- Row 8 ---------------------------------------------------------------
$ws.Cells.Item(8, 8).Value = 2      # H8
$ws.Cells.Item(8, 9).Value = 1      # I8

# --- Row 9 -----------------------------------------------------------------
$ws.Cells.Item(9, 4).Value = 2       # D9
$ws.Cells.Item(9, 8).Value = 5       # H9
$ws.Cells.Item(9, 9).Value = 1       # I9

# --- Row 10 ----------------------------------------------------------------
$ws.Cells.Item(10, 4).Value = 2      # D10
$ws.Cells.Item(10, 8).Value = 5      # H10
$ws.Cells.Item(10, 9).Value = 1      # I10

# --- Row 11 ----------------------------------------------------------------
$ws.Cells.Item(11, 4).Value = 2      # D11
$ws.Cells.Item(11, 8).Value = 10     # H11
$ws.Cells.Item(11, 9).Value = 1      # I11

# --- New row 13 --------------------------------------------------------
# A13 / E13 / G13 hold an explicit empty-text value (not a blank cell) --
# writing a leading apostrophe forces Excel to store a (shared-string)
# empty-text literal instead of clearing the cell, then ClearFormats()
# drops the "quote prefix" display format that the apostrophe trick adds
# so the cell keeps the workbook's default style.
$a13 = $ws.Cells.Item(13, 1)
$a13.Value = "'"
$a13.ClearFormats()

$ws.Cells.Item(13, 2).Value = 21    # B13

$e13 = $ws.Cells.Item(13, 5)
$e13.Value = "'"
$e13.ClearFormats()

$g13 = $ws.Cells.Item(13, 7)
$g13.Value = "'"
$g13.ClearFormats()
